$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new columns E, F, G ---
$ws.Range("E1").Value = "PV diff"
$ws.Range("F1").Value = "Obj corr"
$ws.Range("G1").Value = "Diff"

# --- Fix corrected objective value typo in D4 (9050 -> 9500) ---
$ws.Range("D4").Value = 9500

# --- PV diff column (E) : literal values, 5-decimal number format ---
$ws.Range("E2").Value = 0.40994366040178998
$ws.Range("E3").Value = 0.37253296342661602
$ws.Range("E4").Value = 0.40994369703737998
$ws.Range("E5").Value = 0.37445074093644898
$ws.Range("E6").Value = 0.40994366040198499
$ws.Range("E7").Value = 0.37703065854857298
$ws.Range("E2:E7").NumberFormat = "0.00000"

# --- Obj corr column (F) : C - E, shared formula, 2-decimal number format ---
$ws.Range("F2:F7").Formula = "=C2-E2"
$ws.Range("F2:F7").NumberFormat = "0.00"

# --- Diff column (G) : (F - $F$2), only on rows with the 8-month objective (2,4,6) ---
$ws.Range("G2").Formula = "=(F2-`$F`$2)"
$ws.Range("G4").Formula = "=(F4-`$F`$2)"
$ws.Range("G6").Formula = "=(F6-`$F`$2)"
$ws.Range("G2").Style = "Percent"
$ws.Range("G4").Style = "Percent"
$ws.Range("G6").Style = "Percent"
$ws.Range("G4").NumberFormat = "0.00"
$ws.Range("G6").NumberFormat = "0.00"

# --- H column: relative diff, percentage format ---
$ws.Range("H2").Formula = "=G2/`$F`$2"
$ws.Range("H4").Formula = "=G4/`$F`$2"
$ws.Range("H6").Formula = "=G6/`$F`$2"
$ws.Range("H2").Style = "Percent"
$ws.Range("H4").Style = "Percent"
$ws.Range("H6").Style = "Percent"
